# BOM work, change uart level shifter and m.2 socket
# - Update DDR_A (sheet "DDR_A") flight-time numbers for byte-lane-0 nets
#   (rows 17-19, 21, 23-25): Track Length / Total Length / Track Delay /
#   Total Delay columns (E/G/H/L).
# - Update DDR_B (sheet "DDR_B") rows 19-23: the M.2 socket / UART level
#   shifter net reassignment shifted the net names (column B) down one row
#   (DRAM_D00_B now routes where DRAM_D04_B used to, etc.) and the matching
#   Track Length / Total Length / Track Delay / Package Delay / Total Delay
#   numbers (E/G/H/J/L) moved with them.
# - Active sheet/selection moved from DDR_B (F38) to DDR_A (E41).

$wb = $excel.ActiveWorkbook

$ddrA = $wb.Worksheets.Item("DDR_A")
$ddrB = $wb.Worksheets.Item("DDR_B")

# ---------------------------------------------------------------------
# DDR_A updates (rows 17,18,19,21,23,24,25) - columns E,G,H,L
# ---------------------------------------------------------------------
$ddrAUpdates = @(
    @{ Row = 17; E = 20.79; G = 21.19; H = 141.82; L = 197.82 },
    @{ Row = 18; E = 20.49; G = 20.89; H = 141.53; L = 197.83 },
    @{ Row = 19; E = 19.55; G = 19.95; H = 136.06; L = 197.86 },
    @{ Row = 21; E = 22.23; G = 22.63; H = 150.31; L = 197.91 },
    @{ Row = 23; E = 18.87; G = 19.27; H = 128.51; L = 197.71 },
    @{ Row = 24; E = 21.82; G = 22.22; H = 150.56; L = 197.86 },
    @{ Row = 25; E = 20.07; G = 20.47; H = 138.63; L = 197.83 }
)

foreach ($u in $ddrAUpdates) {
    $r = $u.Row
    $ddrA.Range("E$r").Value = $u.E
    $ddrA.Range("G$r").Value = $u.G
    $ddrA.Range("H$r").Value = $u.H
    $ddrA.Range("L$r").Value = $u.L
}

# ---------------------------------------------------------------------
# DDR_B updates (rows 19-23) - net name (B) + E,G,H,J,L
# ---------------------------------------------------------------------
$ddrBUpdates = @(
    @{ Row = 19; B = "DRAM_D00_B"; E = 16.41; G = 16.41; H = 95.87;  J = 51.5; L = 147.37 },
    @{ Row = 20; B = "DRAM_D04_B"; E = 16.05; G = 16.05; H = 93.73999999999999; J = 53.6; L = 147.34 },
    @{ Row = 21; B = "DRAM_D03_B"; E = 18.04; G = 18.04; H = 105.34; J = 42;   L = 147.34 },
    @{ Row = 22; B = "DRAM_D06_B"; E = 15.84; G = 15.84; H = 92.61;  J = 54.7; L = 147.31 },
    @{ Row = 23; B = "DRAM_D07_B"; E = 16.99; G = 16.99; H = 99.27;  J = 48;   L = 147.27 }
)

foreach ($u in $ddrBUpdates) {
    $r = $u.Row
    $ddrB.Range("B$r").Value = $u.B
    $ddrB.Range("E$r").Value = $u.E
    $ddrB.Range("G$r").Value = $u.G
    $ddrB.Range("H$r").Value = $u.H
    $ddrB.Range("J$r").Value = $u.J
    $ddrB.Range("L$r").Value = $u.L
}

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping: DDR_B's selection moves to L14
# (no longer the active tab) and DDR_A becomes the active tab with its
# selection at E41.
# ---------------------------------------------------------------------
$ddrB.Range("L14").Select()

$ddrA.Activate()
$ddrA.Range("E41").Select()
